# Update "want to go" counts (column F) for three events that appear on
# both the "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet.
#
#   F7  / F9  : 11716 -> 11718
#   F13 / F16 : 5803  -> 5804
#   F15 / F18 : 3525  -> 3526

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value  = 11718
$wsExhibit.Range("F13").Value = 5804
$wsExhibit.Range("F15").Value = 3526

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value  = 11718
$wsAll.Range("F16").Value = 5804
$wsAll.Range("F18").Value = 3526
